$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; this shifts Address/City/State/Zip/Dues
# Paid one column to the right and extends row/column formatting to match.
$ws.Range("B1").EntireColumn.Insert()

# The new column reuses column A's width and formatting (copy format only, no value).
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B1").ColumnWidth = $ws.Range("A1").ColumnWidth

# Add the new "Email" header in row 3, between "Player Name" and "Address".
$ws.Range("B3").Value = "Email"

# Restore the selection to match the authored state.
[void]$ws.Range("B5").Select()
